# Update the "想去人数" (want-to-go count) figures that changed between
# the previous and the newly generated data snapshot.
#
# Sheet "展览" (sheet1):
#   F2  : 120   -> 121
#   F5  : 11156 -> 11161
#   F10 : 11062 -> 11066
#   F16 : 89    -> 90
#
# Sheet "全部类型" (sheet4):
#   F2  : 120   -> 121
#   F7  : 11156 -> 11161
#   F12 : 11062 -> 11066
#   F18 : 89    -> 90

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 121
$wsExhibition.Range("F5").Value = 11161
$wsExhibition.Range("F10").Value = 11066
$wsExhibition.Range("F16").Value = 90

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 121
$wsAll.Range("F7").Value = 11161
$wsAll.Range("F12").Value = 11066
$wsAll.Range("F18").Value = 90
